$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$values = @{
    "B2"  = 0.2066981850894128
    "C2"  = 2.270155390337222
    "D2"  = 13.1891407897574
    "E2"  = 3.631685667807362
    "F2"  = 3.661877248494126
    "G2"  = 51

    "B3"  = -0.075228059080023
    "C3"  = 2.191059243838117
    "D3"  = 12.66358371323103
    "E3"  = 3.558592940086155
    "F3"  = 3.59391839656686
    "G3"  = 50

    "B4"  = 0.207102134124417
    "C4"  = 2.274418628387441
    "D4"  = 13.76964981448061
    "E4"  = 3.710747878053778
    "F4"  = 3.743358472864432
    "G4"  = 49

    "B5"  = 0.03062386419637815
    "C5"  = 2.209507036589459
    "D5"  = 12.80187705469571
    "E5"  = 3.577971080751731
    "F5"  = 3.615701818853944
    "G5"  = 48

    "B6"  = 0.2145225631768136
    "C6"  = 2.299235277238912
    "D6"  = 13.98594729776688
    "E6"  = 3.739779043976646
    "F6"  = 3.77398585940008
    "G6"  = 47

    "B7"  = 0.1255984722002377
    "C7"  = 2.192001147543182
    "D7"  = 12.97327225105472
    "E7"  = 3.601842896498225
    "F7"  = 3.639428754450474
    "G7"  = 46

    "B8"  = 0.1386258507452061
    "C8"  = 2.134562903283757
    "D8"  = 12.94582025312243
    "E8"  = 3.598030051725864
    "F8"  = 3.635985363188866
    "G8"  = 45

    "B9"  = -0.01227215085340897
    "C9"  = 1.819249639236372
    "D9"  = 11.02222735225344
    "E9"  = 3.319973998731532
    "F9"  = 3.358333532090894
    "G9"  = 44

    "B10" = -0.0200094225433071
    "C10" = 1.878919759582091
    "D10" = 11.5868556758706
    "E10" = 3.403947073012534
    "F10" = 3.44417236534218
    "G10" = 43

    "B11" = 0.06995512137203723
    "C11" = 1.589308054450445
    "D11" = 10.76865959017209
    "E11" = 3.28156358923183
    "F11" = 3.320586807489363
    "G11" = 42
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
